# PowerUppXL.xlsx - log new workout entries (22-24 Apr 2019) and refresh
# the aggregate "Exercise Table" summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New dated entries, in the same order the dates were first typed so
#    the shared-string table grows 22/04 -> 23/04 -> 24/04 (as in the
#    original authoring session).
# ---------------------------------------------------------------------

# -- 22/04/2019 entries --------------------------------------------------
$wsReverseLegLift = $wb.Worksheets.Item("Reverse_Leg_Lift")
$wsReverseLegLift.Range("A2").Value = "22/04/2019"
$wsReverseLegLift.Range("B2").Value = 18

$wsSquats = $wb.Worksheets.Item("Squats")
$wsSquats.Range("A2").Value = "22/04/2019"
$wsSquats.Range("B2").Value = 21

# -- 23/04/2019 entries --------------------------------------------------
$wsDumbbellCurls = $wb.Worksheets.Item("Dumbbell_Curls")
$wsDumbbellCurls.Range("A2").Value = "23/04/2019"
$wsDumbbellCurls.Range("B2").Value = 12

$wsDumbbellCurls.Range("A3").Value = "23/04/2019"
$wsDumbbellCurls.Range("B3").Value = 15
$wsDumbbellCurls.Range("B4").Value = 17

# -- 24/04/2019 entries --------------------------------------------------
$wsSitUps = $wb.Worksheets.Item("Sit_Ups")
$wsSitUps.Range("A3").Value = "24/04/2019"
$wsSitUps.Range("B3").Value = 31

$wsPushUps = $wb.Worksheets.Item("Push_Ups")
$wsPushUps.Range("A2").Value = "24/04/2019"

$wsPushUps.Range("B3").Value = 12

$wsPushUps.Range("A4").Value = "24/04/2019"
$wsPushUps.Range("B4").Value = 1

$wsPushUps.Range("A5").Value = "24/04/2019"
$wsPushUps.Range("B5").Value = 12

$wsReverseLegLift.Range("A3").Value = "24/04/2019"
$wsReverseLegLift.Range("B3").Value = 10

$wsSquats.Range("A3").Value = "24/04/2019"
$wsSquats.Range("B3").Value = 36

$wsStandingLunges = $wb.Worksheets.Item("Standing_Lunges")
$wsStandingLunges.Range("A2").Value = "24/04/2019"
$wsStandingLunges.Range("B2").Value = 7
$wsStandingLunges.Range("B3").Value = 10

# -- Dumbbell_Side_Bend keeps a real date serial (17/04/2019), just later --
$wsDumbbellSideBend = $wb.Worksheets.Item("Dumbbell_Side_Bend")
$wsDumbbellSideBend.Range("A2").Value = 43572
$wsDumbbellSideBend.Range("B2").Value = 22

# ---------------------------------------------------------------------
# 2) Refresh the aggregate "Exercise Table" sheet with the new totals.
# ---------------------------------------------------------------------
$wsTable = $wb.Worksheets.Item("Exercise Table")

# Push_Ups row
$wsTable.Range("B2").Value = 99
$wsTable.Range("C2").Value = 99
$wsTable.Range("D2").Value = 12

# Squats row
$wsTable.Range("C3").Value = 99
$wsTable.Range("D3").Value = 36
$wsTable.Range("E3").Value = 21

# Reverse_Leg_Lifts row
$wsTable.Range("C4").Value = 99
$wsTable.Range("D4").Value = 10

# Dumbbell_Side_Bend row
$wsTable.Range("D5").Value = 22

# Dumbbell_Curls row
$wsTable.Range("D6").Value = 17

# Standing_Lunges row
$wsTable.Range("C7").Value = 7
$wsTable.Range("D7").Value = 10

# Sit Ups row
$wsTable.Range("B10").Value = 31

# ---------------------------------------------------------------------
# 3) View-state touch-ups: the Dumbbell_Curls selection moves to its
#    newly-added B4 cell, and Squats becomes the active tab (activated
#    last so it "wins" as the workbook's visible sheet).
# ---------------------------------------------------------------------
$wsDumbbellCurls.Range("B4").Select()

$wsSquats.Activate()
